# MapItem.xlsx - "resource" sheet update
#
# 1. For rows 5-42, column H ("isRuleTile") flips from 1 to 0.
# 2. Columns I/J (rows 5-42) and M/N (rows 5-43) were sharing a style entry
#    that duplicated another existing cellXfs entry; re-applying the
#    (already-in-effect) alignment normalizes each range onto the
#    pre-existing, de-duplicated style so the redundant cellXfs entry is
#    dropped when the workbook is written back out.
# 3. The active selection moves from L41 to H5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 5; $r -le 42; $r++) {
    $ws.Range("H$r").Value = 0

    $iCell = $ws.Range("I$r")
    $iCell.HorizontalAlignment = -4108
    $iCell.VerticalAlignment = -4108

    $jCell = $ws.Range("J$r")
    $jCell.HorizontalAlignment = -4108
    $jCell.VerticalAlignment = -4108
}

for ($r = 5; $r -le 43; $r++) {
    $mCell = $ws.Range("M$r")
    $mCell.HorizontalAlignment = -4108
    $mCell.VerticalAlignment = -4108

    $nCell = $ws.Range("N$r")
    $nCell.HorizontalAlignment = -4108
    $nCell.VerticalAlignment = -4108
}

$ws.Range("H5").Select()
